{"js": "// Remove the trailing \"Ver no Jupiter...\" link paragraph, the\n// \"\u00a9 2020 ... Creative Commons Attribution\" copyright paragraph, and the\n// blank paragraph that immediately followed them (the page footer block\n// that the Jekyll site build no longer emits).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Ver no Jupiter ...\" paragraph by its text so the edit is\n// resilient to any surrounding content shifting line numbers around.\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // The footer block is exactly three consecutive paragraphs:\n  //   1) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n  //   2) \"\u00a9 2020 . Contact: ... Creative Commons Attribution\"\n  //   3) an empty paragraph right after it\n  const toDelete = [];\n  for (let i = startIndex; i < Math.min(startIndex + 3, paragraphs.items.length); i++) {\n    toDelete.push(paragraphs.items[i]);\n  }\n\n  // Delete from the last one back to the first so earlier deletions don't\n  // invalidate the objects we still need to remove.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" link paragraph, the\n# \"(c) 2020 ... Creative Commons Attribution\" copyright paragraph, and the\n# blank paragraph that immediately followed them (the page footer block\n# that the Jekyll site build no longer emits).\n\n$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter ...\" paragraph via Find so the edit is\n# resilient to any surrounding content shifting paragraph indices around.\n$find = $d.Content\n$found = $find.Find.Execute(\"Ver no Jupiter\")\n\nif ($found) {\n    $count = $d.Paragraphs.Count\n    $startIndex = -1\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Start -le $find.Start -and $p.Range.End -ge $find.End) {\n            $startIndex = $i\n            break\n        }\n    }\n\n    if ($startIndex -ne -1) {\n        # The footer block is exactly three consecutive paragraphs:\n        #   1) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n        #   2) \"(c) 2020 . Contact: ... Creative Commons Attribution\"\n        #   3) an empty paragraph right after it\n        $endIndex = [Math]::Min($startIndex + 2, $count)\n        $startRange = $d.Paragraphs.Item($startIndex).Range\n        $endRange = $d.Paragraphs.Item($endIndex).Range\n        $combined = $d.Range($startRange.Start, $endRange.End)\n        $combined.Delete()\n    }\n}\n"}
